$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Extend "Sheet2" (the RGB doubling table) from row 124 down to row 140.
#    Column E continues counting down from 126 to 111; column F = 2*E.
# ---------------------------------------------------------------------------
$rgb = $wb.Worksheets.Item("Sheet2")

for ($r = 125; $r -le 140; $r++) {
    $rgb.Range("E$r").Value = 251 - $r
}

# Rows 125-130 fall inside the pre-existing shared-formula range (F67:F130),
# rows 131-140 are a fresh fill so Excel starts a new shared-formula group.
$rgb.Range("F125:F130").Formula = "=2*E125"
$rgb.Range("F131:F140").Formula = "=2*E131"

# ---------------------------------------------------------------------------
# 2) Add a new sheet "Sheet5" after the last existing sheet. It holds a
#    16-step "wipe" sequence (0..15 in column A) plus a diagonally moving
#    "x" marker used to animate an LED breathe/wipe effect.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wipe = $wb.Worksheets.Add($null, $lastSheet)

$wipeCols = @("I", "H", "G", "F", "E", "D", "C", "B", "B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 1
    $cell = $wipe.Range("A$row")
    $cell.Value = $i
    $cell.VerticalAlignment = -4108
    $wipe.Range($wipeCols[$i] + $row).Value = "x"
}
$wipe.Range("J17").Value = "x"

$wipe.Range("H1").Select() | Out-Null

# The new sheet becomes the active tab.
$wipe.Activate() | Out-Null
